$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 136 (정재헌) is removed; all rows below it shift up by one.
$ws.Rows.Item(136).Delete()
